$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: existing "css" task moves here with its date (45006) already set.
# B6 becomes the (renamed) css task, C6 gets marked "ok".
$ws.Range("B6").Value = "faire le css page d accueil"
$ws.Range("C6").Value = "ok"

# New row 7: auguste's filter-by-name task, dated 45010 (2023-03-25).
# Copy the date formatting from A6 (already using the shared date style)
# instead of assigning NumberFormat directly, so no new style is created.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 45010
$ws.Range("B7").Value = "auguste: filtrer par nom"

# Row 8: yohan's filter-by-ingredient task.
$ws.Range("B8").Value = "yohan: filtrer par ingredient"

# Row 9: jess's filter-by-country task.
$ws.Range("B9").Value = "jess: filtrer par pays"

# Move the active selection to C6, matching the saved view state.
$ws.Range("C6").Select()
